# Update column G ("K") values for rows 2-57 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @(1,0,0,1,0,1,1,0,0,0,1,1,0,0,1,0,0,1,0,2,1,2,1,1,2,0,1,1,0,2,0,0,1,1,1,0,0,0,1,1,0,2,0,0,0,1,1,0,0,1,0,1,1,0,2,2)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
